# Translations check-in for farm_field.xlsx
#
# The ODK display.* settings were renamed to their nested "<kind>.text"
# form (display.text -> display.prompt.text, display.hint -> display.hint.text,
# display.title -> display.title.text, display.new_instance_text ->
# display.new_instance_label.text). On the "survey" sheet the old
# "display.new_instance_text" column (I) was dropped entirely, while on the
# "crops" sheet the equivalent column was kept and just relabeled.
# The "settings" tab also became the active/selected sheet.

$wb = $excel.ActiveWorkbook

# --- "survey" sheet (1st tab) -------------------------------------------
$survey = $wb.Worksheets.Item(1)
$survey.Activate()

# Drop the "display.new_instance_text" column outright; everything to its
# right (required/calculation) shifts one column to the left.
$survey.Columns("I").Delete()

# Rename the remaining display.* headers in place.
$survey.Range("F1").Value = "display.prompt.text"
$survey.Range("G1").Value = "display.hint.text"

$survey.Range("G6").Select()

# --- "model" sheet (2nd tab) ---------------------------------------------
# No cell content changes here.

# --- "settings" sheet (3rd tab) ------------------------------------------
$settings = $wb.Worksheets.Item(3)
$settings.Activate()

$settings.Range("C1").Value = "display.title.text"

$settings.Range("C2").Select()

# --- "queries" sheet (4th tab) --------------------------------------------
# No cell content changes here.

# --- "crops" sheet (5th tab) ----------------------------------------------
$crops = $wb.Worksheets.Item(5)
$crops.Activate()

$crops.Range("F1").Value = "display.prompt.text"
$crops.Range("G1").Value = "display.hint.text"
$crops.Range("J1").Value = "display.new_instance_label.text"

$crops.Range("G3").Select()

# --- "calculates" sheet (6th tab) ------------------------------------------
# No cell content changes here.

# Leave "settings" as the active sheet, matching the saved workbook view.
$settings.Activate()
